$d = $word.ActiveDocument

function Get-MatchRange($para, $text) {
    $rng = $para.Range.Duplicate
    $found = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $text"
    }
    return $rng
}

function Split-At($pos) {
    # Force a clean run split at an exact character position without leaving
    # any bookmark behind: add a temporary bookmark (which splits the run
    # that currently contains $pos) then immediately delete it again.
    $pt = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TmpSplitMarker", $pt) | Out-Null
    $d.Bookmarks("TmpSplitMarker").Delete()
}

# -----------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits right after
#    "Grant Posell".
# -----------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# -----------------------------------------------------------------------
# 2) "Take A Hike" bullet: collapse the split "REST" / "ful" spell-check
#    runs (and drop the now pointless proofErr markers) so the sentence
#    reads "...Google Maps RESTful APIs...".
# -----------------------------------------------------------------------
$hikeParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Google Maps*RESTful*") {
        $hikeParaIndex = $i
        break
    }
}
Write-Host "hike para: $hikeParaIndex"
$hikePara = $d.Paragraphs($hikeParaIndex)

$bootstrapRng = Get-MatchRange $hikePara "Bootstrap, Google Maps "
$bootstrapStart = $bootstrapRng.Start
$afterBootstrap = $bootstrapRng.End

# Rewrite the whole "Bootstrap, Google Maps RESTful APIs...Firebase."
# span in one go. This is the only reliable way to make the engine fully
# reflow the runs and drop the (now stale) spellStart/spellEnd proofErr
# markers that used to bracket "REST"/"ful". Collateral damage: the
# "Bootstrap, Google Maps " segment also gets folded into the resulting
# run, so we re-split it back out below.
$mergeRng = $d.Range($bootstrapStart, $hikePara.Range.End - 1)
Write-Host "merge text before: '$($mergeRng.Text)'"
# Force a genuine text delta first (so the engine actually reflows instead
# of treating an identical re-assignment as a no-op), then restore the
# real wording.
$mergeRng.Text = "Bootstrap, Google Maps RESTfulX APIs, Hike Project API, Open Weather API, and Google FirebaseX."
$mergeRng2 = $d.Range($bootstrapStart, $hikePara.Range.End - 1)
$mergeRng2.Text = "Bootstrap, Google Maps RESTful APIs, Hike Project API, Open Weather API, and Google Firebase."

# Recreate the "Bootstrap, Google Maps " | "REST" | "ful APIs..." run
# boundaries.
Split-At $bootstrapStart
Split-At $afterBootstrap
Split-At ($afterBootstrap + 4)

# -----------------------------------------------------------------------
# 3) "Communicates with clients..." bullet: fix the stray double space
#    and relocate the "_GoBack" bookmark to sit between "...contributes
#    to " and "revenue generating projects.".
# -----------------------------------------------------------------------
$commParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Communicates with clients*") {
        $commParaIndex = $i
        break
    }
}
Write-Host "comm para: $commParaIndex"
$commPara = $d.Paragraphs($commParaIndex)

$commStart = $commPara.Range.Start
$commEnd = $commPara.Range.End - 1
$fullRng = $d.Range($commStart, $commEnd)
Write-Host "comm text before: '$($fullRng.Text)'"
$fullRng.Text = "Communicates with clients to understand their needs and contributes to revenue generating projects."

$theirRng = Get-MatchRange $commPara "their needs"
Split-At $theirRng.Start

$revenueRng = Get-MatchRange $commPara "revenue generating"
$d.Bookmarks.Add("_GoBack", $d.Range($revenueRng.Start, $revenueRng.Start)) | Out-Null

Write-Host "done"
